# Weekly update: insert a new daily price record for Zanahoria (Hortaliza)
# at Terminal La Palmera de La Serena. This shifts the existing data rows
# 83..194 down to 84..195 (preserving their values) and populates the new
# row 83 with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 83 - pushes existing rows 83..194 down to 84..195.
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new record.
$ws.Cells.Item(83, 1).Value = 8
$ws.Cells.Item(83, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(83, 3).Value = "Coquimbo"
$ws.Cells.Item(83, 4).Value = 44467
$ws.Cells.Item(83, 5).Value = 4
$ws.Cells.Item(83, 6).Value = 100114013
$ws.Cells.Item(83, 7).Value = "Zanahoria"
$ws.Cells.Item(83, 8).Value = "Sin especificar"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 600
$ws.Cells.Item(83, 11).Value = 6000
$ws.Cells.Item(83, 12).Value = 7000
$ws.Cells.Item(83, 13).Value = 6500
$ws.Cells.Item(83, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(83, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(83, 16).Value = 325
$ws.Cells.Item(83, 17).Value = 20
$ws.Cells.Item(83, 18).Value = "Hortaliza"

# Make sure the date cell keeps the original date number format used by
# the rest of column D.
$ws.Cells.Item(83, 4).NumberFormat = $ws.Cells.Item(84, 4).NumberFormat
